$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Empresa entry ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Optus"
$ws.Range("C3").Value = 2

# --- Rows 3-4: Usuario entries ---
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2

$ws.Range("F3").Value = "MamazZ"
$ws.Range("F4").Value = "Leleo"

$ws.Range("G3").Value = "5A5"
$ws.Range("G4").Value = "6K7"

$ws.Range("H3").Value = "maria@email.com"
$ws.Range("H4").Value = "leonardo@email.com"

$ws.Range("I3").Value = 123
$ws.Range("I4").Value = 456

$ws.Range("J3").Value = "Comum"
$ws.Range("J4").Value = "Adm"

# --- Rows 3-4: Album entries (idAlbum / idartistas) ---
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 2

$ws.Range("M3").Value = 1
$ws.Range("M4").Value = 2

$ws.Range("N3").Value = 2
$ws.Range("N4").Value = 3

# --- Rows 3-4: EstMusical id refs used by Album ---
$ws.Range("R3").Value = 1
$ws.Range("R4").Value = 2

# --- Rows 8-9: Artistas entries ---
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 2

$ws.Range("F8").Value = "Fefeu"
$ws.Range("F9").Value = "Fufu"

# --- Rows 3-4: Album nomeAlbum ---
$ws.Range("O3").Value = "Solidao"
$ws.Range("O4").Value = "Roll'n'Rock"

# --- Rows 8-10: Cd entries ---
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 2
$ws.Range("L10").Value = 3

$ws.Range("M8").Value = "Nunca Pare"
$ws.Range("M9").Value = "Loucos"
$ws.Range("M10").Value = "Embora..."

$ws.Range("N8").Value = 1
$ws.Range("N9").Value = 2
$ws.Range("N10").Value = 1

$ws.Range("O8").Value = 1
$ws.Range("O9").Value = 3
$ws.Range("O10").Value = 2

# --- Rows 3-5: EstMusical entries ---
$ws.Range("R5").Value = 3

$ws.Range("S5").Value = "Rock"
$ws.Range("S3").Value = "Sertanejo"
$ws.Range("S4").Value = "Modão"

# --- Header fix: C2 now reads "idUsuarios" ---
$ws.Range("C2").Value = "idUsuarios"

# --- Header shuffle on row 2: idCd/idEstiloMusical/nomeAlbum shift one column left ---
$ws.Range("N2").Value = "idEstiloMusical"
$ws.Range("O2").Value = "nomeAlbum"
$ws.Range("P2").ClearContents()

# --- Hyperlinks on the two e-mail cells (Excel auto-hyperlink + "Hiperlink" style) ---
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:maria@email.com")
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:leonardo@email.com")

# --- Selection moved to C2, no frozen/scrolled top-left cell anymore ---
$ws.Range("C2").Select()
